$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 123
$ws.Range("A123").Value = "updel_122"
$ws.Range("B123").Value = "y"
$ws.Range("C123").Value = "主键等值逻辑与删除"
$ws.Range("D123").Value = "SQLFunction"
$ws.Range("F123").Value = "scalar058"
$ws.Range("G123").Value = "scalar_common_value1"
$ws.Range("H123").Value = "delete from `$scalar058 where id=8 and id=28"
$ws.Range("I123").Value = "0"
$ws.Range("J123").Value = "select * from `$scalar058 where id=8 and id=28"
$ws.Range("K123").Value = "src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_122.csv"
$ws.Range("L123").Value = "csv_containsAll"

# Row 124
$ws.Range("A124").Value = "updel_123"
$ws.Range("B124").Value = "y"
$ws.Range("C124").Value = "主键等值逻辑或删除"
$ws.Range("D124").Value = "SQLFunction"
$ws.Range("E124").NumberFormat = "@"
$ws.Range("F124").Value = "scalar058"
$ws.Range("G124").Value = "scalar_common_value1"
$ws.Range("H124").Value = "delete from `$scalar058 where id=8 or id=28"
$ws.Range("I124").Value = "2"
$ws.Range("J124").Value = "select * from `$scalar058 where id=8 or id=28 or id=18 or id=38"
$ws.Range("K124").Value = "src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_123.csv"
$ws.Range("L124").Value = "csv_containsAll"

# Row 125
$ws.Range("A125").Value = "updel_124"
$ws.Range("B125").Value = "y"
$ws.Range("C125").Value = "索引列等值逻辑与删除"
$ws.Range("D125").Value = "SQLFunction"
$ws.Range("F125").Value = "scalar058"
$ws.Range("G125").Value = "scalar_common_value1"
$ws.Range("H125").Value = "delete from `$scalar058 where age=-18 and age=18"
$ws.Range("I125").Value = "0"
$ws.Range("J125").Value = "select age from `$scalar058 where age=-18 and age=18"
$ws.Range("K125").Value = "src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_124.csv"
$ws.Range("L125").Value = "csv_containsAll"

# Row 126
$ws.Range("A126").Value = "updel_125"
$ws.Range("B126").Value = "y"
$ws.Range("C126").Value = "索引列等值逻辑或删除"
$ws.Range("D126").Value = "SQLFunction"
$ws.Range("F126").Value = "scalar058"
$ws.Range("G126").Value = "scalar_common_value1"
$ws.Range("H126").Value = "delete from `$scalar058 where age=-18 or age=18"
$ws.Range("I126").Value = "4"
$ws.Range("J126").Value = "select age from `$scalar058 where age=-18 or age=18"
$ws.Range("K126").Value = "src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_125.csv"
$ws.Range("L126").Value = "csv_containsAll"

# Row 127
$ws.Range("A127").Value = "updel_126"
$ws.Range("B127").Value = "y"
$ws.Range("C127").Value = "主键等值逻辑与更新"
$ws.Range("D127").Value = "SQLFunction"
$ws.Range("F127").Value = "scalar058"
$ws.Range("G127").Value = "scalar_common_value1"
$ws.Range("H127").Value = "update `$scalar058 set age=888 where id=8 and id=28"
$ws.Range("I127").Value = "0"
$ws.Range("J127").Value = "select age from `$scalar058 where id=8 and id=28"
$ws.Range("K127").Value = "src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_126.csv"
$ws.Range("L127").Value = "csv_containsAll"

# Row 128
$ws.Range("A128").Value = "updel_127"
$ws.Range("B128").Value = "y"
$ws.Range("C128").Value = "主键等值逻辑或更新"
$ws.Range("D128").Value = "SQLFunction"
$ws.Range("F128").Value = "scalar058"
$ws.Range("G128").Value = "scalar_common_value1"
$ws.Range("H128").Value = "update `$scalar058 set age=888 where id=8 or id=28"
$ws.Range("I128").Value = "2"
$ws.Range("J128").Value = "select age from `$scalar058 where id=8 or id=28"
$ws.Range("K128").Value = "src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_127.csv"
$ws.Range("L128").Value = "csv_containsAll"

# Row 129
$ws.Range("A129").Value = "updel_128"
$ws.Range("B129").Value = "y"
$ws.Range("C129").Value = "索引列等值逻辑与更新"
$ws.Range("D129").Value = "SQLFunction"
$ws.Range("F129").Value = "scalar058"
$ws.Range("G129").Value = "scalar_common_value1"
$ws.Range("H129").Value = "update `$scalar058 set age=888 where age=-18 and age=18"
$ws.Range("I129").Value = "0"
$ws.Range("J129").Value = "select * from `$scalar058 where age=-18 and age=18"
$ws.Range("K129").Value = "src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_128.csv"
$ws.Range("L129").Value = "csv_containsAll"

# Row 130
$ws.Range("A130").Value = "updel_129"
$ws.Range("B130").Value = "y"
$ws.Range("C130").Value = "索引列等值逻辑或更新"
$ws.Range("D130").Value = "SQLFunction"
$ws.Range("F130").Value = "scalar058"
$ws.Range("G130").Value = "scalar_common_value1"
$ws.Range("H130").Value = "update `$scalar058 set age=888 where age=-18 or age=18"
$ws.Range("I130").Value = "4"
$ws.Range("J130").Value = "select * from `$scalar058 where age=-18 or age=18 or age=888"
$ws.Range("K130").Value = "src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_129.csv"
$ws.Range("L130").Value = "csv_containsAll"

# Column H width: widen to fit the new, longer SQL strings (60.625 -> 88.375 chars)
$ws.Columns.Item(8).ColumnWidth = 87.66071428571429

# Update the active selection to match where the new test cases were entered
$ws.Range("H109").Select()
